# Add input validation: a new participant ("Nawaf Alomari") was missing from
# the "AIoT Hackathon with stc" roster. Insert them as row 8 (team 1 /
# SuperDevops), shifting the two existing team-2 rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AIoT Hackathon with stc")

# Insert a blank row at position 8 (shifts rows 8-9 down to 9-10) without
# disturbing any other columns.
$ws.Range("A8:G8").Insert(-4121)

# Copy the formatting from the row above (row 7) into the freshly inserted row.
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122) | Out-Null

# Fill in the new participant's details.
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 201931050
$ws.Range("C8").Value = "Nawaf Alomari"
$ws.Range("D8").Value = "SWE"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "SuperDevops"

# Column G stores the rank as text ("3.0 "); copy it from an existing cell with
# the same value so it is stored as text rather than being coerced to a number.
$ws.Range("G6").Copy()
$ws.Range("G8").PasteSpecial(-4163) | Out-Null

# Renumber the rows that shifted down.
$ws.Range("A9").Value = 4
$ws.Range("A10").Value = 5
